$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.443.64'
$ws.Range("E2").Value = '  +2.53%  '
$ws.Range("D3").Value = '2.427.00'
$ws.Range("E3").Value = '  +3.19%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '556.87'
$ws.Range("E5").Value = '  +2.36%  '
$ws.Range("D6").Value = '143.83'
$ws.Range("E6").Value = '  +5.10%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  +1.86%  '
$ws.Range("D9").Value = '2.427.77'
$ws.Range("E9").Value = '  +3.30%  '
$ws.Range("E10").Value = '  +5.06%  '
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").Value = '26.37'
$ws.Range("E14").Value = '  +6.69%  '
$ws.Range("E15").Value = '  +9.50%  '
$ws.Range("D16").Value = '2.865.38'
$ws.Range("E16").Value = '  +3.19%  '
$ws.Range("D17").Value = '62.276.40'
$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("D18").Value = '2.427.77'
$ws.Range("E18").Value = '  +3.54%  '
$ws.Range("D19").Value = '11.10'
$ws.Range("E19").Value = '  +4.14%  '
$ws.Range("D20").Value = '324.93'
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("D21").Value = '4.19'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = '6.76'
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '1.79'
$ws.Range("E24").Value = '  +6.60%  '
$ws.Range("D25").Value = '65.00'
$ws.Range("E25").Value = '  +2.58%  '
$ws.Range("D26").Value = '9.09'
$ws.Range("E26").Value = '  +8.78%  '
$ws.Range("D27").Value = '575.21'
$ws.Range("E27").Value = '  +15.21%  '
$ws.Range("D28").Value = '2.546.44'
$ws.Range("E28").Value = '  +3.20%  '
$ws.Range("D29").Value = '0.0₃0948'
$ws.Range("E29").Value = '  +10.20%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +5.85%  '
$ws.Range("E32").Value = '  +6.45%  '
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("E34").Value = '  +4.16%  '
$ws.Range("E35").Value = '  +5.78%  '
$ws.Range("D36").Value = '5.74'
$ws.Range("E36").Value = '  +9.47%  '
$ws.Range("E37").Value = '  +5.71%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").Value = '0.385'
$ws.Range("E39").Value = '  +2.60%  '
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("E41").Value = '  +1.85%  '
$ws.Range("D42").Value = '150.55'
$ws.Range("E42").Value = '  +5.04%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '41.71'
$ws.Range("E44").Value = '  +2.74%  '
$ws.Range("E45").Value = '  +15.21%  '
$ws.Range("D46").Value = '151.18'
$ws.Range("E46").Value = '  +5.55%  '
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("E48").Value = '  +4.90%  '
$ws.Range("D49").Value = '20.50'
$ws.Range("E49").Value = '  +7.46%  '
$ws.Range("E50").Value = '  +3.91%  '
$ws.Range("D51").Value = '0.0919'
$ws.Range("E51").Value = '  +2.19%  '
